$wb = $excel.ActiveWorkbook

# Remember the originally active sheet so we can restore it at the end
# (selecting a range on another sheet activates that sheet as a side effect).
$originalActive = $wb.ActiveSheet.Name

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Cells.Item(5, 1).Value = 0.97
$wsSummary.Cells.Item(5, 5).Value = 0.97
$wsSummary.Cells.Item(5, 6).Value = 0.97
$wsSummary.Range("A6:F6").Clear()
$wsSummary.Cells.Item(2, 7).Value = 5
$wsSummary.Cells.Item(2, 7).Value = ""
$wsSummary.Cells.Item(2, 7).Style = "Normal"
$wsSummary.Range("D5").Select()

# --- Repayment schedule sheet ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# P2 removed entirely
$wsRepay.Cells.Item(2, 16).Clear()

$wsRepay.Cells.Item(3, 10).Value = 0.51
$wsRepay.Cells.Item(3, 11).Value = 888.23
# O3 removed entirely
$wsRepay.Cells.Item(3, 15).Clear()
$wsRepay.Cells.Item(3, 16).Value = 888.23

$wsRepay.Cells.Item(4, 10).Value = 0.46
$wsRepay.Cells.Item(4, 11).Value = 888.18
# O4 removed entirely
$wsRepay.Cells.Item(4, 15).Clear()
$wsRepay.Cells.Item(4, 16).Value = 888.18

$wsRepay.Cells.Item(5, 10).Value = 0
$wsRepay.Cells.Item(5, 11).Value = 887.72
# O5 removed entirely
$wsRepay.Cells.Item(5, 15).Clear()
$wsRepay.Cells.Item(5, 16).Value = 887.72

# O6, O7, O8 removed entirely
$wsRepay.Cells.Item(6, 15).Clear()
$wsRepay.Cells.Item(7, 15).Clear()
$wsRepay.Cells.Item(8, 15).Clear()

$wsRepay.Range("J8").Select()

# --- Transactions sheet ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("G2").Select()

# Restore the originally active sheet
$wb.Worksheets.Item($originalActive).Activate()
